$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column D (this shifts the existing D:K financial-year
#    columns one place right, to E:L), ready to receive a new "2018" column
#    of figures at the front of each table (Income Statement, Balance Sheet,
#    Cash Flow Statement).
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).Insert()

# ---------------------------------------------------------------------------
# 2) "Period Ending" header rows: new column D gets the date 2018-12-31
#    (serial 43465), formatted the same as the existing date cells.
# ---------------------------------------------------------------------------
$dateRows = @(7, 38, 80)
foreach ($r in $dateRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = "[$-409]d\-mmm\-yy;@"
    $cell.Value = 43465
}

# ---------------------------------------------------------------------------
# 3) Blank separator rows: column D still needs the same numeric style as
#    the rest of the row even though it stays empty.
# ---------------------------------------------------------------------------
$blankRows = @(11, 16, 19, 39, 40, 55, 56, 67, 82, 90, 95)
foreach ($r in $blankRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = "#,##0"
    $cell.HorizontalAlignment = -4152
}

# ---------------------------------------------------------------------------
# 4) New figures for column D (the 2018 financial year) across the three
#    statements.
# ---------------------------------------------------------------------------
$newYearData = @{
    8 = 22500
    9 = 7500
    10 = 15000
    12 = 7700
    13 = 0
    14 = 0
    15 = 0
    17 = 38100
    18 = -15600
    20 = -2200
    21 = -17400
    22 = 0
    23 = -17800
    24 = -600
    25 = 0
    26 = -17200
    27 = -17200
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 2200
    33 = -17200
    34 = 0
    35 = -17200
    41 = 22400
    42 = "NA"
    43 = 3900
    44 = 800
    45 = 1100
    46 = 28300
    47 = 0
    48 = 1700
    49 = 2600
    50 = 0
    51 = 0
    52 = 100
    53 = 0
    54 = 32700
    57 = 1500
    58 = 700
    59 = 4400
    60 = 6500
    61 = 9300
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 15800
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -169500
    73 = 0
    74 = 0
    75 = 0
    76 = 16900
    77 = 0
    81 = -17200
    83 = 400
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = -10800
    91 = -700
    92 = 0
    93 = 0
    94 = -1500
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 17500
    101 = -100
    102 = 5000
}

foreach ($r in $newYearData.Keys) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = "#,##0"
    $cell.HorizontalAlignment = -4152
    $cell.Value = $newYearData[$r]
}

# ---------------------------------------------------------------------------
# 5) A couple of historical "Capital Expenditures" figures (row 91) were
#    corrected at the same time, not just shifted over from the old layout.
# ---------------------------------------------------------------------------
$ws.Cells.Item(91, 5).Value = -1000
$ws.Cells.Item(91, 6).Value = -100
